$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-31 Wednesday" "2026-01-01 Thursday"

Replace-Text "625×4=2500" "769×5=3845"
Replace-Text "269×9=2421" "691×9=6219"
Replace-Text "682×4=2728" "351×9=3159"
Replace-Text "113×4=452" "241×7=1687"
Replace-Text "703×6=4218" "356×7=2492"

Replace-Text "517×9=4653" "528×2=1056"
Replace-Text "230×8=1840" "269×5=1345"
Replace-Text "666×7=4662" "275×7=1925"
Replace-Text "468×7=3276" "778×3=2334"
Replace-Text "538×4=2152" "117×2=234"

Replace-Text "821×7=5747" "356×8=2848"
Replace-Text "894×8=7152" "872×3=2616"
Replace-Text "505×5=2525" "768×3=2304"
Replace-Text "564×9=5076" "491×9=4419"
Replace-Text "280×9=2520" "541×7=3787"

Replace-Text "115×2=230" "902×3=2706"
Replace-Text "205×3=615" "599×8=4792"
Replace-Text "973×8=7784" "834×3=2502"
Replace-Text "842×5=4210" "326×6=1956"
Replace-Text "622×6=3732" "672×7=4704"

Replace-Text "812×3=2436" "456×6=2736"
Replace-Text "259×2=518" "668×6=4008"
Replace-Text "133×5=665" "813×4=3252"
Replace-Text "482×7=3374" "445×7=3115"
Replace-Text "435×5=2175" "740×4=2960"
